# Updates the cryptos list with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we touch to stay as text, since many of the
# values (e.g. "1.003") would otherwise be auto-coerced into numbers by Excel.
$priceCells = "D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D39","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51"
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.433.30"
$ws.Range("E2").Value = "  -3.77%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.769.91"
$ws.Range("E3").Value = "  -2.91%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5 - USDC
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.04%  "

# Row 6 - BNB
$ws.Range("D6").Value = "306.23"
$ws.Range("E6").Value = "  -2.07%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4275"
$ws.Range("E7").Value = "  +0.76%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3676"
$ws.Range("E8").Value = "  +1.88%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.07182"
$ws.Range("E9").Value = "  -0.25%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "0.8476"
$ws.Range("E10").Value = "  -1.38%  "

# Row 11 - Solana
$ws.Range("D11").Value = "20.32"
$ws.Range("E11").Value = "  -1.38%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.772.17"
$ws.Range("E12").Value = "  -3.16%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "5.252"
$ws.Range("E13").Value = "  -2.60%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "6.437"
$ws.Range("E14").Value = "  -0.52%  "

# Row 15 - TRON
$ws.Range("D15").Value = "0.06815"
$ws.Range("E15").Value = "  -1.72%  "

# Row 16 - BinanceUSD
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.17%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "79.42"
$ws.Range("E17").Value = "  -1.17%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.000008646"
$ws.Range("E18").Value = "  -2.36%  "

# Row 19 - Dai
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.07%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "15.02"
$ws.Range("E20").Value = "  -2.21%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "26.431.17"
$ws.Range("E21").Value = "  -3.85%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "5.086"
$ws.Range("E22").Value = "  -0.58%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "11.25"
$ws.Range("E23").Value = "  +3.05%  "

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range("D24").Value = "2.003.54"
$ws.Range("E24").Value = "  -2.41%  "

# Row 25 - Monero
$ws.Range("D25").Value = "152.51"
$ws.Range("E25").Value = "  -1.66%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "1.850"
$ws.Range("E26").Value = "  -6.77%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "18.14"
$ws.Range("E27").Value = "  -2.91%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "5.095"
$ws.Range("E28").Value = "  -0.90%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "114.46"
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 - LidoDAOToken
$ws.Range("D30").Value = "1.708"
$ws.Range("E30").Value = "  -5.01%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "0.08938"
$ws.Range("E31").Value = "  +0.95%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "0.7278"
$ws.Range("E32").Value = "  -2.35%  "

# Row 33 - was ARBITRUM, now Filecoin (rows 33/34 swapped order + new data)
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "4.345"
$ws.Range("E33").Value = "  -4.07%  "

# Row 34 - was Filecoin, now ARBITRUM
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.114"
$ws.Range("E34").Value = "  -0.66%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "2.757"
$ws.Range("E35").Value = "  -7.33%  "

# Row 36 - Frax
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  +0.05%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").Value = "1.077"
$ws.Range("E37").Value = "  -0.69%  "

# Row 38 - Hedera (Volume only)
$ws.Range("E38").Value = "  -2.61%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "0.01897"
$ws.Range("E39").Value = "  -1.34%  "

# Row 40 - Algorand (Volume only)
$ws.Range("E40").Value = "  -2.00%  "

# Row 41 - TheSandbox
$ws.Range("D41").Value = "0.4922"
$ws.Range("E41").Value = "  -2.68%  "

# Row 42 - MXToken
$ws.Range("D42").Value = "2.535"
$ws.Range("E42").Value = "  -9.21%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "6.222"
$ws.Range("E43").Value = "  -3.29%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "8.061"
$ws.Range("E44").Value = "  -3.17%  "

# Row 45 - Quant
$ws.Range("D45").Value = "104.90"
$ws.Range("E45").Value = "  -0.99%  "

# Row 46 - was PaxDollar, now EnergySwap (rows 46/47 swapped order + new data)
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.20"
$ws.Range("E46").Value = "  -2.58%  "

# Row 47 - was EnergySwap, now PaxDollar
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  +0.07%  "

# Row 48 - Cronos
$ws.Range("D48").Value = "0.06193"
$ws.Range("E48").Value = "  -4.07%  "

# Row 49 - Decentraland
$ws.Range("D49").Value = "0.4483"
$ws.Range("E49").Value = "  -3.86%  "

# Row 50 - NEARProtocol
$ws.Range("D50").Value = "1.583"
$ws.Range("E50").Value = "  -1.59%  "

# Row 51 - RenderToken
$ws.Range("D51").Value = "1.750"
$ws.Range("E51").Value = "  +3.55%  "
